$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 687; this shifts existing rows 687-728 down to 688-729
# and extends the sheet dimension from A1:D728 to A1:D729.
$ws.Rows.Item(687).Insert()

# Column A holds the date as plain text (e.g. "2026/12/29"), not a real Date value.
# Setting .Value directly on a date-like string auto-converts it into a serial date,
# so force a text number format first, then restore the plain (unstyled) look of the
# surrounding data cells afterwards.
$newDateCell = $ws.Range("A687")
$newDateCell.NumberFormat = "@"
$newDateCell.Value = "2026/01/20"
$newDateCell.Style = $ws.Range("A688").Style

$ws.Range("B687").Value = "火"
$ws.Range("C687").Value = 9
$ws.Range("D687").Value = 143
